$d = $word.ActiveDocument

# Locate the "CU011: " label (collapses $r to the matched range on success).
$r = $d.Content
$found = $r.Find.Execute("CU011: ", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($found) {
    $insStart = $r.End
    $newText = "Ver notificación"

    # Insert the new text right after "CU011: ", before the paragraph mark.
    $ins = $d.Range($insStart, $insStart)
    $ins.InsertAfter($newText)

    # Force Word to materialize the inserted text as its own run (rather
    # than being silently coalesced into the preceding "CU011: " run) by
    # toggling a character property on it and then reverting the toggle,
    # which leaves the run's effective formatting unchanged (bCs, sz 24,
    # szCs 24 inherited from the surrounding text).
    $newRange = $d.Range($insStart, $insStart + $newText.Length)
    $newRange.Font.Bold = $true
    $newRange.Font.Bold = $false
}
